# Updates the cryptocurrency price/volume(1h) table on Sheet1 to reflect
# the latest scrape (GitHub Actions scheduled refresh). Coin names/links in
# B:C only change where two rows swapped rank position; D (Price) and
# E (Volume(1h)) are refreshed for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds free-form text (e.g. "69.136.64", "0.0₃0741") rather
# than numbers, same as the rest of the sheet. Force text formatting first so
# Excel's autodetect doesn't coerce numeric-looking values like "579.82" or
# "0.120" into floating-point numbers and lose the exact printed digits.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.136.64"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.387.16"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "579.82"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "178.54"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +8.08%  "
$ws.Range("D10").Value = "0.585"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "48.17"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").Value = "682.86"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "8.58"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.930.68"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "69.214.00"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.377.73"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "17.65"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "11.27"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "0.907"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "101.06"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "9.70"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").Value = "33.42"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "8.72"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "3.81"
$ws.Range("E31").Value = "  +15.90%  "
$ws.Range("D32").Value = "11.02"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "552.27"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "57.81"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "3.601.43"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "35.37"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "0.0₃0741"
$ws.Range("E40").Value = "  +10.52%  "
$ws.Range("E41").Value = "  +4.41%  "
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  +4.78%  "
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "130.30"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "2.60"
$ws.Range("E51").Value = "  +1.92%  "
